$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2..347) holds a date serial ("Förändrad" = last changed date).
# Bump every value in that range from 45177 to 45178 (one day later),
# mirroring an automated daily refresh of the sheet.
$lastRow = 347
$rng = $ws.Range("C2:C$lastRow")
$rng.Value = 45178
